# This edit moves the "Senior Analyst - Myers Research" work-experience block so
# that it appears immediately after "Partner - Siege Analytics" (i.e. before the
# "Research Director - PCCC" block), and moves the "Field Director - The Feldman
# Group" block to fill the position the Myers Research block vacated (i.e.
# immediately after the "Research Director - PCCC" block, before "Programmer -
# Lake Research Partners").
#
# All four job blocks (PCCC, Myers Research, Lake Research Partners, Feldman
# Group) share an identical paragraph shape: one Heading3 title paragraph
# followed by one Normal subtitle paragraph and three Normal bullet
# paragraphs. Because no paragraphs are being added or removed overall, this
# reorder is implemented by rewriting the run text of paragraphs 15-34 (the
# four contiguous job blocks following "Partner - Siege Analytics") so they
# read, in order: Myers Research, PCCC, Feldman Group, Lake Research Partners.

$d = $word.ActiveDocument

$newText = @{
  15 = "Senior Analyst - Myers Research (Austin, TX) | 2012 - 2014"
  16 = "Political Research & Analysis"
  17 = [char]0x2022 + " Designed comprehensive survey instruments for specialized voting segments and niche markets"
  18 = [char]0x2022 + " Developed sophisticated analytical products and reports that delivered actionable insights to clients"
  19 = [char]0x2022 + " Co-developed a web application to manage all aspects of survey operations, from instrument design to data collection and analysis"

  20 = "Research Director - PCCC (Washington, DC) | August 2011 - August 2012"
  21 = "Political Research & Data Analysis (FLEEM System)"
  22 = [char]0x2022 + " Conceived, architected, and engineered FLEEM web application using Twilio API handling tens of thousands of simultaneous phone calls using emulated predictive dialer for regulated political surveys"
  23 = [char]0x2022 + " Developed IVR polling system for early quantitative research supporting Senators Martin Heinrich and Elizabeth Warren"
  24 = [char]0x2022 + " Built comprehensive tabular and graphical reporting system with Python, GeoDjango, PostGIS, and Apache webserver"

  25 = "Field Director - The Feldman Group (Austin, TX) | 2011 - 2012"
  26 = "Political Campaign Management"
  27 = [char]0x2022 + " Managed all aspects of survey fielding for a multi-million dollar research firm, including scheduling, oversight, sampling, and quality control"
  28 = [char]0x2022 + " Developed and implemented data warehousing solutions for efficient storage and retrieval of research findings"
  29 = [char]0x2022 + " Created custom reports and data visualizations based on specific client requirements"

  30 = "Programmer - Lake Research Partners (Austin, TX) | 2008"
  31 = "Political Polling & Research"
  32 = [char]0x2022 + " Designed questionnaires and analyzed data for complex market research studies across diverse industries"
  33 = [char]0x2022 + " Conducted statistical modeling and analysis to address multifaceted consumer behavior questions"
  34 = [char]0x2022 + " Pioneered the integration of advanced mapping techniques into standard reports, including choropleths and hexagonal grid maps"
}

foreach ($i in 15..34) {
    $para = $d.Paragraphs.Item($i)
    $rng = $para.Range
    # Keep the paragraph mark, only replace the visible text portion.
    $rng.MoveEnd(1, -1) | Out-Null
    $rng.Text = $newText[$i]
}

Write-Host "Reordered professional-experience blocks (Myers Research / PCCC / Feldman Group / Lake Research Partners)."
